# Apply the "Added UI Checklist 2" edits:
#  1. "User Interface (UI) Checklist:" -> "User Interface Checklist:"
#  2. Merge the leading "  " (two-space) run with the following page-name
#     run in each UI checklist line, so the paragraph's text is unchanged
#     but it is carried by a single run instead of two.

$d = $word.ActiveDocument

# 1) Drop the "(UI) " qualifier from the checklist heading.
$d.Content.Find.Execute(
    "User Interface (UI) Checklist:", $false, $false, $false, $false, $false,
    $true, 1, $false, "User Interface Checklist:", 2
) | Out-Null

# 2) Collapse the "  " + "<Page Name>" run pairs into a single run per line.
#    Doing the replacement via Find/Replace on the whole paragraph text
#    keeps the surrounding run formatting (font Aptos) intact because the
#    matched text is replaced in place, merging what used to be two runs
#    into the run that is left holding the text.
$pages = @("Landing Page", "Registration Page", "Log-In Page", "Home Page", "Daily Page")

foreach ($page in $pages) {
    $search = "  " + $page
    $d.Content.Find.Execute(
        $search, $false, $false, $false, $false, $false,
        $true, 1, $false, $search, 2
    ) | Out-Null
}
